$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.040790813259761
$ws.Range("D2").Value = 1.04916559464119
$ws.Range("E2").Value = 0.992614727750844
$ws.Range("F2").Value = 1.05775577049272
$ws.Range("I2").Value = 1.042924056688473
$ws.Range("J2").Value = 1.045875514491193
$ws.Range("K2").Value = 1.051923411866803
$ws.Range("L2").Value = 0.9955398523335997
$ws.Range("M2").Value = 1.060489870676634
$ws.Range("N2").Value = 1.047360777685368

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.041608523233425
$ws.Range("D3").Value = 1.0498152761327
$ws.Range("E3").Value = 0.9936372048519299
$ws.Range("F3").Value = 1.058535772695139
$ws.Range("I3").Value = 1.043116535728473
$ws.Range("J3").Value = 1.046339644603778
$ws.Range("K3").Value = 1.052385704830444
$ws.Range("L3").Value = 0.9963617723202687
$ws.Range("M3").Value = 1.061083857760616
$ws.Range("N3").Value = 1.04782556691595

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.042138205475054
$ws.Range("D4").Value = 1.050236177153582
$ws.Range("E4").Value = 0.9942998659930998
$ws.Range("F4").Value = 1.0590413480743
$ws.Range("I4").Value = 1.043240223967725
$ws.Range("J4").Value = 1.046639854253768
$ws.Range("K4").Value = 1.052684665220367
$ws.Range("L4").Value = 0.9968940712668347
$ws.Range("M4").Value = 1.061468406863188
$ws.Range("N4").Value = 1.048126202898078

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.042361018451689
$ws.Range("D5").Value = 1.050413245085541
$ws.Range("E5").Value = 0.994578699834602
$ws.Range("F5").Value = 1.059254096282758
$ws.Range("I5").Value = 1.043292016220359
$ws.Range("J5").Value = 1.046766033922337
$ws.Range("K5").Value = 1.052810305140397
$ws.Range("L5").Value = 0.9971179600053012
$ws.Range("M5").Value = 1.061630117493185
$ws.Range("N5").Value = 1.048252561756249

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.042398437572085
$ws.Range("D6").Value = 1.05044298262117
$ws.Range("E6").Value = 0.994625531979634
$ws.Range("F6").Value = 1.059289829566381
$ws.Range("I6").Value = 1.043300700249456
$ws.Range("J6").Value = 1.046787218337307
$ws.Range("K6").Value = 1.052831398073232
$ws.Range("L6").Value = 0.9971555583673455
$ws.Range("M6").Value = 1.061657272080536
$ws.Range("N6").Value = 1.048273776255518

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.04214118218635
$ws.Range("D7").Value = 1.050238542670244
$ws.Range("E7").Value = 0.994303590798249
$ws.Range("F7").Value = 1.059044190025948
$ws.Range("I7").Value = 1.043240916829813
$ws.Range("J7").Value = 1.04664154038614
$ws.Range("K7").Value = 1.052686344197898
$ws.Range("L7").Value = 0.9968970624462089
$ws.Range("M7").Value = 1.061470567469104
$ws.Range("N7").Value = 1.048127891424951

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.041067043287728
$ws.Range("D8").Value = 1.049385050090079
$ws.Range("E8").Value = 0.9929600610674297
$ws.Range("F8").Value = 1.058019196247965
$ws.Range("I8").Value = 1.042989283208738
$ws.Range("J8").Value = 1.046032392181397
$ws.Range("K8").Value = 1.052079681166008
$ws.Range("L8").Value = 0.9958175282591056
$ws.Range("M8").Value = 1.060690569054661
$ws.Range("N8").Value = 1.047517878159886

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.039178708366584
$ws.Range("D9").Value = 1.04788510181819
$ws.Range("E9").Value = 0.9906006454969559
$ws.Range("F9").Value = 1.056219714880493
$ws.Range("I9").Value = 1.042539329757926
$ws.Range("J9").Value = 1.044958182012588
$ws.Range("K9").Value = 1.051009390146103
$ws.Range("L9").Value = 0.9939188001724441
$ws.Range("M9").Value = 1.059317714084664
$ws.Range("N9").Value = 1.046442142489425

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.037922902085167
$ws.Range("D10").Value = 1.046887944581202
$ws.Range("E10").Value = 0.989033133672735
$ws.Range("F10").Value = 1.055024676907355
$ws.Range("I10").Value = 1.042235005947414
$ws.Range("J10").Value = 1.044241571464587
$ws.Range("K10").Value = 1.050295087592437
$ws.Range("L10").Value = 0.9926553831429383
$ws.Range("M10").Value = 1.05840365180001
$ws.Range("N10").Value = 1.045724514272249

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.037379876929471
$ws.Range("D11").Value = 1.046456853630861
$ws.Range("E11").Value = 0.988355674866747
$ws.Range("F11").Value = 1.054508332399292
$ws.Range("I11").Value = 1.042102206786333
$ws.Range("J11").Value = 1.04393117570164
$ws.Range("K11").Value = 1.049985619900519
$ws.Range("L11").Value = 0.9921088820399291
$ws.Range("M11").Value = 1.058008151846565
$ws.Range("N11").Value = 1.045413677711716

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.037178287195223
$ws.Range("D12").Value = 1.046296831997426
$ws.Range("E12").Value = 0.9881042295826724
$ws.Range("F12").Value = 1.054316708593475
$ws.Range("I12").Value = 1.042052726017027
$ws.Range("J12").Value = 1.043815867364552
$ws.Range("K12").Value = 1.049870645657122
$ws.Range("L12").Value = 0.9919059725120875
$ws.Range("M12").Value = 1.057861291542284
$ws.Range("N12").Value = 1.045298205623563

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.037221523691737
$ws.Range("D13").Value = 1.046331152411317
$ws.Range("E13").Value = 0.9881581567098651
$ws.Range("F13").Value = 1.054357804845675
$ws.Range("I13").Value = 1.042063346732937
$ws.Range("J13").Value = 1.043840601985795
$ws.Range("K13").Value = 1.049895309094034
$ws.Range("L13").Value = 0.9919494934313052
$ws.Range("M13").Value = 1.057892791468372
$ws.Range("N13").Value = 1.045322975370805

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.037363211123895
$ws.Range("D14").Value = 1.046443624045259
$ws.Range("E14").Value = 0.9883348863814464
$ws.Range("F14").Value = 1.054492489245444
$ws.Range("I14").Value = 1.042098119812444
$ws.Range("J14").Value = 1.043921644547062
$ws.Range("K14").Value = 1.049976116578684
$ws.Range("L14").Value = 0.9920921077337197
$ws.Range("M14").Value = 1.057996011380257
$ws.Range("N14").Value = 1.045404133021805

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.037450524582994
$ws.Range("D15").Value = 1.046512935465882
$ws.Range("E15").Value = 0.9884438009545853
$ws.Range("F15").Value = 1.054575495275442
$ws.Range("I15").Value = 1.042119524369645
$ws.Range("J15").Value = 1.04397157578958
$ws.Range("K15").Value = 1.050025901576114
$ws.Range("L15").Value = 0.9921799884222134
$ws.Range("M15").Value = 1.058059614718176
$ws.Range("N15").Value = 1.045454135172415

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.037958956766886
$ws.Range("D16").Value = 1.046916569253566
$ws.Range("E16").Value = 0.9890781214508737
$ws.Range("F16").Value = 1.055058968654442
$ws.Range("I16").Value = 1.042243797858967
$ws.Range("J16").Value = 1.044262169440701
$ws.Range("K16").Value = 1.050315622480012
$ws.Range("L16").Value = 0.9926916645766087
$ws.Range("M16").Value = 1.058429906187442
$ws.Range("N16").Value = 1.045745141499852

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.038278084117506
$ws.Range("D17").Value = 1.047169942702385
$ws.Range("E17").Value = 0.989476357848556
$ws.Range("F17").Value = 1.055362538645458
$ws.Range("I17").Value = 1.042321477449748
$ws.Range("J17").Value = 1.044444425657669
$ws.Range("K17").Value = 1.050497312211481
$ws.Range("L17").Value = 0.9930127773699352
$ws.Range("M17").Value = 1.058662260408191
$ws.Range("N17").Value = 1.045927656541552

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.038464297685291
$ws.Range("D18").Value = 1.047317797088036
$ws.Range("E18").Value = 0.9897087662937556
$ws.Range("F18").Value = 1.055539713312124
$ws.Range("I18").Value = 1.042366687699101
$ws.Range("J18").Value = 1.044550723018858
$ws.Range("K18").Value = 1.050603272265966
$ws.Range("L18").Value = 0.9932001317071769
$ws.Range("M18").Value = 1.05879781705054
$ws.Range("N18").Value = 1.046034104857187

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.038527803857033
$ws.Range("D19").Value = 1.047368222740252
$ws.Range("E19").Value = 0.9897880325774034
$ws.Range("F19").Value = 1.055600143454065
$ws.Range("I19").Value = 1.042382086421131
$ws.Range("J19").Value = 1.044586966003258
$ws.Range("K19").Value = 1.050639399067204
$ws.Range("L19").Value = 0.9932640239640975
$ws.Range("M19").Value = 1.058844043121712
$ws.Range("N19").Value = 1.046070399310782

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.038243837314981
$ws.Range("D20").Value = 1.047142751295291
$ws.Range("E20").Value = 0.9894336180360679
$ws.Range("F20").Value = 1.055329957331956
$ws.Range("I20").Value = 1.042313153393494
$ws.Range("J20").Value = 1.044424872274772
$ws.Range("K20").Value = 1.050477820318723
$ws.Range("L20").Value = 0.9929783193494215
$ws.Range("M20").Value = 1.058637328054259
$ws.Range("N20").Value = 1.04590807539061

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.037321484546166
$ws.Range("D21").Value = 1.046410501046012
$ws.Range("E21").Value = 0.9882828385668249
$ws.Range("F21").Value = 1.054452823350225
$ws.Range("I21").Value = 1.042087884231685
$ws.Range("J21").Value = 1.04389777989582
$ws.Range("K21").Value = 1.049952321444537
$ws.Range("L21").Value = 0.9920501090198102
$ws.Range("M21").Value = 1.057965614406069
$ws.Range("N21").Value = 1.045380234480021

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.036742225098356
$ws.Range("D22").Value = 1.045950713167374
$ws.Range("E22").Value = 0.9875604150241495
$ws.Range("F22").Value = 1.053902316289596
$ws.Range("I22").Value = 1.041945362289667
$ws.Range("J22").Value = 1.043566298354654
$ws.Range("K22").Value = 1.049621780321764
$ws.Range("L22").Value = 0.9914670000341481
$ws.Range("M22").Value = 1.057543548269412
$ws.Range("N22").Value = 1.045048282197045

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.037049238323142
$ws.Range("D23").Value = 1.046194397269119
$ws.Range("E23").Value = 0.9879432794643023
$ws.Range("F23").Value = 1.054194056744186
$ws.Range("I23").Value = 1.042020999629852
$ws.Range("J23").Value = 1.043742029919214
$ws.Range("K23").Value = 1.049797019153948
$ws.Range("L23").Value = 0.991776070289318
$ws.Range("M23").Value = 1.057767267661444
$ws.Range("N23").Value = 1.045224263320583

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.0382593117487
$ws.Range("D24").Value = 1.047155037718218
$ws.Range("E24").Value = 0.9894529299347244
$ws.Range("F24").Value = 1.055344679093418
$ws.Range("I24").Value = 1.042316914982016
$ws.Range("J24").Value = 1.044433707636683
$ws.Range("K24").Value = 1.050486627917173
$ws.Range("L24").Value = 0.9929938892766442
$ws.Range("M24").Value = 1.058648593824397
$ws.Range("N24").Value = 1.045916923299748

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.039666352151153
$ws.Range("D25").Value = 1.048272387395707
$ws.Range("E25").Value = 0.9912096547607049
$ws.Range("F25").Value = 1.056684119021552
$ws.Range("I25").Value = 1.042656424324025
$ws.Range("J25").Value = 1.045235979595359
$ws.Range("K25").Value = 1.051286228583985
$ws.Range("L25").Value = 0.9944092447426414
$ws.Range("M25").Value = 1.059672430136972
$ws.Range("N25").Value = 1.046720334576626
